$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '64.640.51'
$ws.Range('E2').Value = '  -0.58%  '
$ws.Range('D3').Value = '3.141.44'
$ws.Range('E3').Value = '  +1.02%  '
$ws.Range('E4').Value = '  +0.16%  '
$ws.Range('D5').Value = '571.96'
$ws.Range('E5').Value = '  +0.78%  '
$ws.Range('D6').Value = '147.92'
$ws.Range('E6').Value = '  -0.57%  '
$ws.Range('E7').Value = '  +0.03%  '
$ws.Range('D8').Value = '3.142.48'
$ws.Range('E8').Value = '  +1.28%  '
$ws.Range('D9').Value = '0.523'
$ws.Range('E9').Value = '  +0.08%  '
$ws.Range('D10').Value = '0.157'
$ws.Range('E10').Value = '  -2.48%  '
$ws.Range('D11').Value = '6.04'
$ws.Range('E11').Value = '  -1.57%  '
$ws.Range('E12').Value = '  +0.04%  '
$ws.Range('D13').Value = '0.0000257'
$ws.Range('E13').Value = '  +3.65%  '
$ws.Range('D14').Value = '36.76'
$ws.Range('E14').Value = '  -0.32%  '
$ws.Range('D15').Value = '3.654.17'
$ws.Range('E15').Value = '  +1.21%  '
$ws.Range('D16').Value = '64.825.54'
$ws.Range('E16').Value = '  -0.01%  '
$ws.Range('D17').Value = '3.139.18'
$ws.Range('E17').Value = '  +1.22%  '
$ws.Range('D18').Value = '7.04'
$ws.Range('E18').Value = '  -0.62%  '
$ws.Range('E19').Value = '  -0.14%  '
$ws.Range('D20').Value = '498.23'
$ws.Range('E20').Value = '  -0.84%  '
$ws.Range('D21').Value = '14.69'
$ws.Range('E21').Value = '  -0.18%  '
$ws.Range('D22').Value = '0.708'
$ws.Range('E22').Value = '  -0.94%  '
$ws.Range('D23').Value = '15.11'
$ws.Range('E23').Value = '  -3.26%  '
$ws.Range('E24').Value = '  -1.54%  '
$ws.Range('D25').Value = '83.52'
$ws.Range('E25').Value = '  -1.16%  '
$ws.Range('E26').Value = '  -0.19%  '
$ws.Range('D27').Value = '8.81'
$ws.Range('E27').Value = '  +1.36%  '
$ws.Range('E28').Value = '  -0.70%  '
$ws.Range('E29').Value = '  +1.11%  '
$ws.Range('E30').Value = '  +2.92%  '
$ws.Range('D31').Value = '27.30'
$ws.Range('E31').Value = '  -1.28%  '
$ws.Range('D32').Value = '0.998'
$ws.Range('E32').Value = '  +0.04%  '
$ws.Range('E33').Value = '  -0.01%  '
$ws.Range('E34').Value = '  +2.02%  '
$ws.Range('E35').Value = '  -1.92%  '
$ws.Range('D36').Value = '54.27'
$ws.Range('E36').Value = '  -2.28%  '
$ws.Range('D37').Value = '0.0891'
$ws.Range('E37').Value = '  +5.97%  '
$ws.Range('D38').Value = '463.11'
$ws.Range('E38').Value = '  +0.01%  '
$ws.Range('D39').Value = '0.0412'
$ws.Range('E39').Value = '  -1.81%  '
$ws.Range('D40').Value = '2.94'
$ws.Range('E40').Value = '  +2.93%  '
$ws.Range('E41').Value = '  +0.78%  '
$ws.Range('D42').Value = '3.016.03'
$ws.Range('E42').Value = '  -2.79%  '
$ws.Range('D43').Value = '0.115'
$ws.Range('E43').Value = '  -4.45%  '
$ws.Range('B44').Value = 'TheGraph'
$ws.Range('C44').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range('D44').Value = '0.281'
$ws.Range('E44').Value = '  -1.62%  '
$ws.Range('B45').Value = 'Fetch.AI'
$ws.Range('C45').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D45').Value = '2.41'
$ws.Range('E45').Value = '  +0.91%  '
$ws.Range('D46').Value = '28.03'
$ws.Range('E46').Value = '  -2.63%  '
$ws.Range('D47').Value = '0.0₃0570'
$ws.Range('E47').Value = '  +4.30%  '
$ws.Range('E49').Value = '  -1.43%  '
$ws.Range('E50').Value = '  -0.81%  '
$ws.Range('D51').Value = '117.81'
$ws.Range('E51').Value = '  -0.15%  '
